# Update Excel files after daily scrape - 2025-08-18 03:37:46 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adjust column widths that changed. The ColumnWidth property (character
# units) serializes to the sheet XML's "width" attribute with a constant
# +5/6 padding offset, so the requested values are pre-compensated here to
# land exactly on the target XML widths of 57 / 89 / 17 / 16 / 40.
$ws.Columns.Item(3).ColumnWidth = 56.166666666666664
$ws.Columns.Item(4).ColumnWidth = 88.16666666666667
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 39.166666666666664

# New data for rows 2-10 (A:H)
$data = @(
    @("1326946", "https://aiesec.org/opportunity/global-talent/1326946", "[Impact Brazil] - Recruitment and Selection Specialist", "São Paulo, SP, Brasil", "No", "0 applicants", "6 - 18 Months", "Hiring"),
    @("1326944", "https://aiesec.org/opportunity/global-talent/1326944", "Digital Marketing Executive", "Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "TIE innovated Solution"),
    @("1326917", "https://aiesec.org/opportunity/global-talent/1326917", "Business Development Intern", "Ahmedabad, Gujarat, India", "No", "1 applicant", "3 - 6 Months", "Port The Container Pvt. Ltd."),
    @("1324592", "https://aiesec.org/opportunity/global-talent/1324592", "Digital marketing", "New Damietta City, Damietta El-Gadeeda City, New Damietta, Damietta Governorate, Egypt", "No", "13 applicants", "3 - 6 Months", "Business Haven Consultancy"),
    @("1322882", "https://aiesec.org/opportunity/global-talent/1322882", "BUSINESS DEVELOPMENT", "Gebze, Kocaeli, Türkiye", "No", "97 applicants", "3 - 6 Months", "ÖZLER KALIP VE İSKELE SİSTEMLERİ A.Ş."),
    @("1313206", "https://aiesec.org/opportunity/global-talent/1313206", "Digital Media Strategist", "Colombo, Sri Lanka", "No", "41 applicants", "9 - 12 Weeks", "Brand Corridor (Pvt) Ltd"),
    @("1310446", "https://aiesec.org/opportunity/global-talent/1310446", "Education Coordinator", "Bursa, Türkiye", "No", "29 applicants", "9 - 12 Weeks", "Genç Kardelen Kindergarden"),
    @("1307741", "https://aiesec.org/opportunity/global-talent/1307741", "Marketing Intern", "Cyberjaya, Selangor, Malaysia", "No", "177 applicants", "6 - 18 Months", "IX Telecom Sdn Bhd"),
    @("1305878", "https://aiesec.org/opportunity/global-talent/1305878", "Sales and Marketing Intern", "Petaling Jaya, Selangor, Malaysia", "No", "156 applicants", "3 - 6 Months", "Business Media International")
)

$rowIndex = 2
foreach ($row in $data) {
    # Column A holds opportunity IDs that must stay text (not auto-converted to
    # a number), so force text entry with a leading apostrophe, then reset the
    # cell style back to Normal so no stray number-format style lingers.
    $ws.Cells.Item($rowIndex, 1).Value = "'" + $row[0]
    $ws.Cells.Item($rowIndex, 1).Style = "Normal"
    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
